$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has rows 53-55 describing "Black Amber" lots
# (Especial/Primera/Segunda) sold "$/caja 15 kilos empedrada" dated 44214.
# The edit adds three NEW rows for the same lots but sold as
# "$/bandeja 10 kilos granel" dated 44571, inserted right before the
# existing row 56, pushing all the following rows down by three.

# Step 1: insert three blank rows at position 56 (rows 53-55 remain intact,
# rows 56 onward shift down by three).
$ws.Rows("56:58").Insert()

# Step 2: the original rows 53:55 content needs to reappear at 56:58 (since
# they were pushed down conceptually). Copy the still-intact 53:55 block
# into the newly inserted blank rows 56:58.
$ws.Range("A53:T55").Copy()
$ws.Range("A56").PasteSpecial()
$excel.CutCopyMode = $false

# Step 3: overwrite rows 53:55 in place with the new "bandeja" records.
# Row 53: Black Amber / Especial
$ws.Range("D53").Value2 = 44571
$ws.Range("M53").Value2 = 65
$ws.Range("Q53").Value2 = "$/bandeja 10 kilos granel"
$ws.Range("S53").Value2 = 1300
$ws.Range("T53").Value2 = 10

# Row 54: Black Amber / Primera
$ws.Range("D54").Value2 = 44571
$ws.Range("M54").Value2 = 60
$ws.Range("Q54").Value2 = "$/bandeja 10 kilos granel"
$ws.Range("S54").Value2 = 1200
$ws.Range("T54").Value2 = 10

# Row 55: Black Amber / Segunda
$ws.Range("D55").Value2 = 44571
$ws.Range("M55").Value2 = 67
$ws.Range("Q55").Value2 = "$/bandeja 10 kilos granel"
$ws.Range("S55").Value2 = 1000
$ws.Range("T55").Value2 = 10
